# Replace the three M2Doc field-code paragraphs (fldChar begin/instrText/fldChar end)
# with plain-text paragraphs using the "{m:...}" token syntax, as produced by the
# TokenIteratorFieldRewriterSplit parser.
#
# Each target paragraph currently holds a single Word field whose code is the
# M2Doc instruction (" m:for p | Sequence{...}", " m:p ", " m:endfor "). We
# rebuild each paragraph from scratch as a sequence of <w:t> runs (mirroring the
# original instrText run-splitting) wrapped in "{" / "}" token delimiters, and
# use Range.InsertXML to swap the whole paragraph content (field included) for
# the new plain-text runs.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-RunXml([string]$text, [bool]$preserve) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    if ($preserve) {
        return "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    } else {
        return "<w:r><w:t>$escaped</w:t></w:r>"
    }
}

function New-ParaXml($runs) {
    $inner = ""
    foreach ($run in $runs) {
        $inner += New-RunXml $run.Text $run.Preserve
    }
    return "<w:p $wNs>$inner</w:p>"
}

# --- Paragraph holding the "for" field -------------------------------------
$forRuns = @(
    @{ Text = "{m:"; Preserve = $false }
    @{ Text = "for p | "; Preserve = $true }
    @{ Text = "Sequence{'Page1'"; Preserve = $false }
    @{ Text = ", 'newPage'.asPagination()"; Preserve = $false }
    @{ Text = ", 'Page2'"; Preserve = $false }
    @{ Text = ", 'newPage'.asPagination()"; Preserve = $false }
    @{ Text = ", 'Page3'"; Preserve = $false }
    @{ Text = ", 'newPage'.asPagination()"; Preserve = $false }
    @{ Text = ", 'Page4'}}"; Preserve = $true }
)

# --- Paragraph holding the "m:p" field --------------------------------------
$pRuns = @(
    @{ Text = "{m:p}"; Preserve = $true }
)

# --- Paragraph holding the "endfor" field -----------------------------------
$endforRuns = @(
    @{ Text = "{m:endfor}"; Preserve = $true }
)

$targets = @(
    @{ Match = "for p |"; Runs = $forRuns },
    @{ Match = "m:p";     Runs = $pRuns },
    @{ Match = "endfor";  Runs = $endforRuns }
)

# $Paragraph.Range.Fields is not reliable for locating the field owned by a
# given paragraph in this host, so resolve the field -> paragraph mapping
# through the document's Fields collection (which does report correct Code
# text per field) and match it back to the paragraph at the same 1-based
# position (fields appear in the same document order as their paragraphs
# here: field 1 in paragraph 2, field 2 in paragraph 3, field 3 in
# paragraph 4).
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $f = $d.Fields.Item($i)
    $code = $f.Code.Text

    foreach ($target in $targets) {
        if ($code.Contains($target.Match)) {
            # Paragraph index = field index + 1 (paragraph 1 is the intro
            # text before any field).
            $p = $d.Paragraphs.Item($i + 1)
            $xml = New-ParaXml $target.Runs
            $p.Range.InsertXML($xml)
            break
        }
    }
}
